$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I20").Value = "ba"
$ws.Range("J20").Value = "Appreciation"
$ws.Range("I22").Value = "sv"
$ws.Range("J22").Value = "Statement-opinion"
$ws.Range("I41").Value = "%"
$ws.Range("J41").Value = "Uninterpretable"
$ws.Range("I57").Value = "sd"
$ws.Range("J57").Value = "Statement-non-opinion"
$ws.Range("I93").Value = "ba"
$ws.Range("J93").Value = "Appreciation"
$ws.Range("I104").Value = "sd"
$ws.Range("J104").Value = "Statement-non-opinion"
$ws.Range("I120").Value = "sv"
$ws.Range("J120").Value = "Statement-opinion"
$ws.Range("I123").Value = "sv"
$ws.Range("J123").Value = "Statement-opinion"
$ws.Range("I127").Value = "sd"
$ws.Range("J127").Value = "Statement-non-opinion"
$ws.Range("I128").Value = "sd"
$ws.Range("J128").Value = "Statement-non-opinion"
$ws.Range("I165").Value = "sd"
$ws.Range("J165").Value = "Statement-non-opinion"
$ws.Range("I192").Value = "b"
$ws.Range("J192").Value = "Acknowledge (Backchannel)"
$ws.Range("I216").Value = "sd"
$ws.Range("J216").Value = "Statement-non-opinion"
$ws.Range("I233").Value = "%"
$ws.Range("J233").Value = "Uninterpretable"
$ws.Range("I235").Value = "%"
$ws.Range("J235").Value = "Uninterpretable"
$ws.Range("I264").Value = "ba"
$ws.Range("J264").Value = "Appreciation"
$ws.Range("I268").Value = "aa"
$ws.Range("J268").Value = "Agree/Accept"
$ws.Range("I271").Value = "ba"
$ws.Range("J271").Value = "Appreciation"
$ws.Range("I299").Value = "sv"
$ws.Range("J299").Value = "Statement-opinion"
$ws.Range("I300").Value = "sd"
$ws.Range("J300").Value = "Statement-non-opinion"
$ws.Range("I305").Value = "sv"
$ws.Range("J305").Value = "Statement-opinion"
$ws.Range("I324").Value = "b"
$ws.Range("J324").Value = "Acknowledge (Backchannel)"
$ws.Range("I329").Value = "b"
$ws.Range("J329").Value = "Acknowledge (Backchannel)"
$ws.Range("I341").Value = "qy"
$ws.Range("J341").Value = "Yes-No-Question"
$ws.Range("I371").Value = "b"
$ws.Range("J371").Value = "Acknowledge (Backchannel)"
$ws.Range("I372").Value = "b"
$ws.Range("J372").Value = "Acknowledge (Backchannel)"
$ws.Range("I379").Value = "b"
$ws.Range("J379").Value = "Acknowledge (Backchannel)"
$ws.Range("I387").Value = "sd"
$ws.Range("J387").Value = "Statement-non-opinion"
$ws.Range("I390").Value = "b"
$ws.Range("J390").Value = "Acknowledge (Backchannel)"
$ws.Range("I393").Value = "sv"
$ws.Range("J393").Value = "Statement-opinion"
$ws.Range("I396").Value = "sd"
$ws.Range("J396").Value = "Statement-non-opinion"
$ws.Range("I398").Value = "sv"
$ws.Range("J398").Value = "Statement-opinion"
$ws.Range("I403").Value = "ba"
$ws.Range("J403").Value = "Appreciation"
$ws.Range("I412").Value = "%"
$ws.Range("J412").Value = "Uninterpretable"
$ws.Range("I420").Value = "b"
$ws.Range("J420").Value = "Acknowledge (Backchannel)"
$ws.Range("I421").Value = "sd"
$ws.Range("J421").Value = "Statement-non-opinion"
$ws.Range("I432").Value = "ba"
$ws.Range("J432").Value = "Appreciation"
$ws.Range("I438").Value = "sd"
$ws.Range("J438").Value = "Statement-non-opinion"
$ws.Range("I462").Value = "b"
$ws.Range("J462").Value = "Acknowledge (Backchannel)"
